$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM line for the 1uF capacitor ("C4,C15,C17,C28") is being split into
# two lines: C15 is a DNP (do-not-place) part that should get its own row,
# separate from the populated C4, C17, C28 parts.

# Insert a new row below row 4 (the existing "C4,C15,C17,C28" row) to hold
# the new C15 line; this shifts rows 5.. down by one.
$ws.Rows("5:5").Insert()

# New row 5 holds the split-off C15 part: same Value/Footprint/LCSC#/Part#
# as row 4, quantity 1, and marked DNP.
$ws.Range("A5").Value = "C15"
$ws.Range("B5").Value = "1uF"
$ws.Range("C5").Value = "Capacitor_SMD:C_0805_2012Metric_Pad1.18x1.45mm_HandSolder"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "DNP"
$ws.Range("F5").Value = "C91185"
$ws.Range("G5").Value = "CC0805KKX7R9BB105"

# Row 4 keeps the same Value/Footprint/LCSC#/Part# but now only references
# C4, C17, C28 with a quantity of 3 (C15 moved out).
$ws.Range("A4").Value = "C4,C17,C28"
$ws.Range("D4").Value = 3

# Selection in the saved file lands on H5 (matches the authored workbook).
$ws.Range("H5").Select()
